$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = @(5660,5620,5580,5540,5520,5510,5500,5490,5500,5510,5520,5540,5560,5590,5620,5660,5720,5790,5880,6000,6140,6290,6460,6640,6820,7000,7160,7300,7410,7490,7540,7560,7540,7500,7430,7340,7240,7130,7020,6920,6820,6730,6640,6570,6500,6440,6390,6360,6330,6310,6300,6310,6320,6340,6370,6410,6450,6500,6560,6630,6710,6810,6910,7020,7140,7250,7360,7470,7570,7680,7780,7880,7970,8040,8070,8060,8030,7980,7910,7800,7670,7550,7420,7300,7170,7030,6900,6750,6620,6480,6340,6220,6020,5950,5920,5860)

for ($i = 0; $i -lt $newA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newA[$i]
    $oldB = $ws.Cells.Item($row, 2).Value2
    $ws.Cells.Item($row, 2).Value = $oldB + 1
}

"done"
